# Roteiros.xlsx — reorganize analysis folders, regenerate test 1/2/3 data,
# and update per-route map selections.

$wb = $excel.ActiveWorkbook

# --- "1 Entregador": move selection from E10 to F12 ---
$ws1 = $wb.Worksheets.Item("1 Entregador")
$ws1.Activate()
$ws1.Range("F12").Select()

# --- "2 Entregadores": move selection from E10 to I3 ---
$ws2 = $wb.Worksheets.Item("2 Entregadores")
$ws2.Activate()
$ws2.Range("I3").Select()

# --- "3 Entregadores": regenerated descarga (E column) test data,
#     becomes the active/selected tab, selection moves to E14 ---
$ws3 = $wb.Worksheets.Item("3 Entregadores")
$ws3.Activate()

$ws3.Range("E3").Value = 220
$ws3.Range("E4").Value = 140
$ws3.Range("E5").Value = 120
$ws3.Range("E6").Value = 210
$ws3.Range("E10").Value = 300
$ws3.Range("E12").Value = 110
$ws3.Range("E13").Value = 70

$ws3.Range("E14").Select()

# --- "4 Entregadores": move selection from E5 to G14 ---
$ws4 = $wb.Worksheets.Item("4 Entregadores")
$ws4.Activate()
$ws4.Range("G14").Select()

# --- "Não é realizar a rota": select entire sheet (all rows/cols) ---
$ws5 = $wb.Worksheets.Item("Não é realizar a rota")
$ws5.Activate()
$ws5.Range("A1:XFD1048576").Select()

# --- "Não é possível com apenas 1": no longer the active tab (selection unchanged) ---
$ws6 = $wb.Worksheets.Item("Não é possível com apenas 1")
$ws6.Activate()

# Final active tab is "3 Entregadores" (activeTab points at it, tabSelected="1")
$ws3.Activate()
